$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns for zh-cn (B2) and de-de (C2), plus the
#     de-de "Latest Handoff Date" (D2) which shares the same timestamp as the
#     de-de sheet's "Latest Handoff Datetime" cell.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-03-23 07:55:14"

# --- zh-cn sheet: Status (C2) and Latest Handoff Datetime (E2)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-23 07:55:06"

# --- de-de sheet: Status (C2) and Latest Handoff Datetime (E2)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-23 07:55:14"
